$wb = $excel.ActiveWorkbook

# The "participants" sheet had a row for NIMS scan 102 (nims_title
# "20170608_15227") that is no longer wanted - remove it entirely. The two
# rows below it (scans 103 and 104) shift up to take rows 4 and 5.
$participants = $wb.Worksheets.Item("participants")
$participants.Rows.Item(4).Delete()

# Make "participants" the active tab/sheet (it was "protocol" before), and
# leave a specific selection on it, matching where the user's cursor ended
# up after doing this edit.
$participants.Activate()
$participants.Range("C15").Select()
